$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 79; this shifts the existing rows 79-159 down to 80-160
# (mirrors the diff, where a new weekly record is prepended and everything below
# moves down by one row).
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new weekly record.
$ws.Range("A79").Value2 = 4
$ws.Range("B79").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C79").Value = "Los Lagos"
$ws.Range("D79").Value2 = 44512
$ws.Range("E79").Value2 = 10
$ws.Range("F79").Value = "Fruta"
$ws.Range("G79").Value2 = 100108
$ws.Range("H79").Value = "Tropicales y subtropicales"
$ws.Range("I79").Value2 = 100108005
$ws.Range("J79").Value = "Piña"
$ws.Range("K79").Value = "Caramelo"
$ws.Range("L79").Value = "Tercera"
$ws.Range("M79").Value2 = 200
$ws.Range("N79").Value2 = 22000
$ws.Range("O79").Value2 = 23000
$ws.Range("P79").Value2 = 22500
$ws.Range("Q79").Value = "$/caja 16 unidades"
$ws.Range("R79").Value = "Ecuador"
$ws.Range("S79").Value2 = 1406
$ws.Range("T79").Value2 = 16
